$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.986.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'207.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "'0.489"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'22.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "'1.782.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "'1.545.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "'3.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "'62.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'26.981.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("D19").Value = "'217.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "'7.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").Value = "'9.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  -1.82%  "
$ws.Range("D25").Value = "'153.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "'15.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("D29").Value = "'1.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").Value = "'1.424.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("E34").Value = "  +3.84%  "
$ws.Range("E35").Value = "  +3.23%  "
$ws.Range("D36").Value = "'1.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.75%  "
$ws.Range("E37").Value = "  +1.30%  "
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("E39").Value = "  +1.67%  "
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("D45").Value = "'64.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").Value = "'1.695.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").Value = "'87.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("D49").Value = "'0.0522"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "'0.0₆0100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  -0.10%  "
